$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "artur"
$ws.Range("B5").Value = "krause"
$ws.Range("C5").Value = "123456789-00"
$ws.Range("D5").Value = "artur@gmail.com"

# E5 is text "123" (not a number) in the target workbook, so force a text
# number format before assigning, then restore the default cell style so no
# stray style index is left behind on the cell.
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "123"
$ws.Range("E5").Style = "Normal"
